$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 299, shifting existing rows 299..375 down to 300..376
$ws.Rows.Item(299).Insert()

# Populate the new row 299 with the new data
$ws.Cells.Item(299, 1).Value = 11
$ws.Cells.Item(299, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(299, 3).Value = "Bíobío"
$ws.Cells.Item(299, 4).Value = 44932
$ws.Cells.Item(299, 5).Value = 8
$ws.Cells.Item(299, 6).Value = 100114001
$ws.Cells.Item(299, 7).Value = "Papa"
$ws.Cells.Item(299, 8).Value = "Asterix"
$ws.Cells.Item(299, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(299, 10).Value = 350
$ws.Cells.Item(299, 11).Value = 13000
$ws.Cells.Item(299, 12).Value = 14000
$ws.Cells.Item(299, 13).Value = 13429
$ws.Cells.Item(299, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(299, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(299, 16).Value = 537
$ws.Cells.Item(299, 17).Value = 25
$ws.Cells.Item(299, 18).Value = "Hortaliza"
